$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the refreshed cryptos feed.
# Numeric-looking price strings (e.g. "254.71") are prefixed with a leading
# apostrophe so Excel keeps storing them as text, matching the original
# inline-string cell contents instead of auto-converting them to numbers.
$ws.Range("D2").Value = "42.723.83"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.211.74"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'254.71"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").Value = "'0.609"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").Value = "'75.36"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("D10").Value = "'40.93"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "'6.88"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "2.534.88"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'14.27"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "2.204.64"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "42.639.21"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "'0.0000103"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "'71.06"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'5.94"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "'2.20"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'227.94"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "'9.40"
$ws.Range("E24").Value = "  -8.16%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'10.56"
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'39.15"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "'173.13"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'20.21"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'0.0844"
$ws.Range("E33").Value = "  +6.26%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").Value = "'0.0348"
$ws.Range("E37").Value = "  +5.45%  "
$ws.Range("D38").Value = "'4.31"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "'12.37"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "  +18.83%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'60.24"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").Value = "'5.26"
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("D44").Value = "'0.197"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("E45").Value = "  -4.99%  "
$ws.Range("D46").Value = "'8.36"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "'0.464"
$ws.Range("E48").Value = "  +5.03%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.436.85"
$ws.Range("E51").Value = "  -0.27%  "
